$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates to the "Price" (column D) and "Volume(1h)" (column E) columns.
# `ForceText` marks Price values that look like plain numbers to Excel's
# auto-detection (e.g. "1.030", "0.07210") so trailing zeros / exponent
# rewriting don't corrupt the text-formatted price strings; values such as
# "27.615.78" already contain two dots and Excel keeps them as text natively.
$updates = @(
    @{ Row = 2; D = "27.615.78"; E = "  +2.62%  "; ForceText = $false },
    @{ Row = 3; D = "1.854.23"; E = $null; ForceText = $false },
    @{ Row = 4; D = "1.030"; E = "  +2.54%  "; ForceText = $true },
    @{ Row = 5; D = "321.85"; E = "  +3.32%  "; ForceText = $true },
    @{ Row = 6; D = $null; E = "  +2.53%  "; ForceText = $false },
    @{ Row = 7; D = "0.4397"; E = "  +2.56%  "; ForceText = $true },
    @{ Row = 8; D = "0.3793"; E = "  +3.08%  "; ForceText = $true },
    @{ Row = 9; D = "0.07429"; E = "  +2.73%  "; ForceText = $true },
    @{ Row = 10; D = "0.8797"; E = "  +2.20%  "; ForceText = $true },
    @{ Row = 11; D = "21.69"; E = "  +2.71%  "; ForceText = $true },
    @{ Row = 12; D = "1.867.49"; E = "  -8.51%  "; ForceText = $false },
    @{ Row = 13; D = "5.533"; E = "  +2.83%  "; ForceText = $true },
    @{ Row = 14; D = "6.709"; E = "  +1.24%  "; ForceText = $true },
    @{ Row = 15; D = "0.07210"; E = "  +4.68%  "; ForceText = $true },
    @{ Row = 16; D = "83.23"; E = "  +3.21%  "; ForceText = $true },
    @{ Row = 17; D = $null; E = "  +3.07%  "; ForceText = $false },
    @{ Row = 18; D = "0.000009063"; E = "  +2.39%  "; ForceText = $true },
    @{ Row = 19; D = $null; E = "  +2.50%  "; ForceText = $false },
    @{ Row = 20; D = "15.46"; E = "  +1.88%  "; ForceText = $true },
    @{ Row = 21; D = "27.648.66"; E = "  +2.57%  "; ForceText = $false },
    @{ Row = 22; D = "5.277"; E = "  +1.78%  "; ForceText = $true },
    @{ Row = 23; D = "11.40"; E = "  +3.69%  "; ForceText = $true },
    @{ Row = 24; D = "158.19"; E = "  +2.87%  "; ForceText = $true },
    @{ Row = 25; D = "1.920"; E = "  +1.89%  "; ForceText = $true },
    @{ Row = 26; D = "18.78"; E = "  +2.82%  "; ForceText = $true },
    @{ Row = 27; D = "1.984"; E = "  +5.24%  "; ForceText = $true },
    @{ Row = 28; D = "5.301"; E = "  +1.62%  "; ForceText = $true },
    @{ Row = 29; D = "117.22"; E = "  +2.04%  "; ForceText = $true },
    @{ Row = 30; D = $null; E = "  +1.51%  "; ForceText = $false },
    @{ Row = 31; D = "1.206"; E = "  +4.17%  "; ForceText = $true },
    @{ Row = 32; D = "0.7642"; E = "  +2.83%  "; ForceText = $true },
    @{ Row = 33; D = $null; E = "  +2.71%  "; ForceText = $false },
    @{ Row = 34; D = $null; E = "  +3.01%  "; ForceText = $false },
    @{ Row = 35; D = $null; E = "  +2.10%  "; ForceText = $false },
    @{ Row = 36; D = "1.153"; E = "  +3.30%  "; ForceText = $true },
    @{ Row = 37; D = "0.01983"; E = "  +3.33%  "; ForceText = $true },
    @{ Row = 38; D = "0.05321"; E = "  +2.09%  "; ForceText = $true },
    @{ Row = 39; D = "0.5174"; E = "  +1.81%  "; ForceText = $true },
    @{ Row = 40; D = "2.818"; E = "  +1.76%  "; ForceText = $true },
    @{ Row = 41; D = "0.1680"; E = "  +2.29%  "; ForceText = $true },
    @{ Row = 42; D = "6.801"; E = "  +5.76%  "; ForceText = $true },
    @{ Row = 43; D = "8.530"; E = "  +3.08%  "; ForceText = $true },
    @{ Row = 44; D = "109.04"; E = "  +2.06%  "; ForceText = $true },
    @{ Row = 45; D = "10.58"; E = "  +1.87%  "; ForceText = $true },
    @{ Row = 46; D = "1.716"; E = "  +4.00%  "; ForceText = $true },
    @{ Row = 47; D = "0.4661"; E = "  +2.30%  "; ForceText = $true },
    @{ Row = 48; D = "0.06401"; E = "  +1.83%  "; ForceText = $true },
    @{ Row = 49; D = "1.860"; E = "  +3.11%  "; ForceText = $true },
    @{ Row = 50; D = "39.42"; E = "  +4.36%  "; ForceText = $true },
    @{ Row = 51; D = "64.17"; E = "  +0.71%  "; ForceText = $true }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)
        if ($u.ForceText) {
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
